$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Mark items 12-14 (rows 15-17) as Done ---
$ws.Range("C15").Value = "Done"
$ws.Range("C16").Value = "Done"
$ws.Range("C17").Value = "Done"

# --- 2. Append new outstanding tasks (IDs 16-24, rows 19-27) ---
$newTasks = @(
    "Disable delete/save buttons when doing DB operations",
    "Add wait cursor on all delete/save operations",
    "Add default 'manual tax' to DB configuration row",
    "Add default 'debt age' to DB configuration row",
    "Add default 'supplier' to invoice screen (add this default to the DB configuration row)",
    "Setup confluence on ezeeit.com server and add basic instructions for ezee invoices",
    "Add 'invoice date' to invoice entity (= date on the supplier invoice)",
    "Add invoice date from/to search filter fields on the invoice grid",
    "Ass invoice date to the grid model"
)

$startRow = 19
$startId = 16

for ($i = 0; $i -lt $newTasks.Length; $i++) {
    $row = $startRow + $i
    $id = $startId + $i

    $ws.Cells.Item($row, 1).Value = $id
    $ws.Cells.Item($row, 2).Value = $newTasks[$i]
    $ws.Cells.Item($row, 3).Value = "Not Done"

    $rowRange = $ws.Range("A$row`:C$row")
    $rowRange.Borders.Weight = -4138
    $ws.Rows.Item($row).RowHeight = 16
}

# --- 3. Extend conditional formatting ranges to cover the new rows ---
$fcs = $ws.Range("C4").FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fcs.Item($i).ModifyAppliesToRange($ws.Range("C4:C27"))
}

Write-Output "edits applied"
